$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

for ($r = 1; $r -le 259; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value()
    if ($v -eq $oldVal) {
        $cell.Value = $newVal
    }
}
